$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview!G3 and de-de!H3 share the same text value (generation timestamp for
# the 3f49a6c0... handback entry) and both move from 14:54:57 -> 14:56:10.
$wsOverview.Range("G3").Value = "2016-08-13 14:56:10"
$wsDeDe.Range("H3").Value = "2016-08-13 14:56:10"

# zh-cn!H3 (Correspond Handoff Datetime) and zh-cn!K3 (Correspond Handback
# DateTime) for the same row get refreshed timestamps.
$wsZhCn.Range("H3").Value = "2016-08-13 14:56:01"
$wsZhCn.Range("K3").Value = "2016-08-13 14:56:26"

# de-de!K3 (Correspond Handback DateTime) refreshed timestamp.
$wsDeDe.Range("K3").Value = "2016-08-13 14:56:36"
